{"js": "// Update the date line and every math-problem cell in the table to the\n// new values from the commit (document-order 1:1 text replacement).\n\nconst DATE_NEW = \"2025-10-05 Sunday\";\nconst TABLE_NEW = [\n  [\"34+23=\", \"55-35=\", \"16+31=\", \"52-4=\", \"99-34=\"],\n  [\"33-1=\", \"7+21=\", \"79-35=\", \"53+30=\", \"31-19=\"],\n  [\"1+35=\", \"73-17=\", \"97-72=\", \"46-22=\", \"23-6=\"],\n  [\"2+35=\", \"45+40=\", \"25-13=\", \"67-25=\", \"43+33=\"],\n  [\"2+40=\", \"67+9=\", \"68-34=\", \"46-16=\", \"55-6=\"],\n  [\"43+25=\", \"56-9=\", \"69+27=\", \"99-7=\", \"51+23=\"],\n  [\"45-4=\", \"1+46=\", \"30+33=\", \"99-73=\", \"3+60=\"],\n  [\"89-19=\", \"5+43=\", \"97-47=\", \"26-23=\", \"86-22=\"],\n  [\"4+54=\", \"51-5=\", \"81-81=\", \"87-21=\", \"92-70=\"],\n  [\"88-27=\", \"12+38=\", \"93-55=\", \"86-20=\", \"92-15=\"],\n  [\"41+35=\", \"97-42=\", \"53-16=\", \"31-8=\", \"10+82=\"],\n  [\"47-23=\", \"19-15=\", \"84-56=\", \"53-9=\", \"15+42=\"],\n  [\"72-13=\", \"98-1=\", \"78+17=\", \"92-79=\", \"18+70=\"],\n  [\"81-25=\", \"73-20=\", \"21+77=\", \"11+36=\", \"81-47=\"],\n  [\"78+18=\", \"46-9=\", \"30+42=\", \"3+72=\", \"43+46=\"],\n  [\"7+56=\", \"35-7=\", \"24+67=\", \"92-89=\", \"40+50=\"],\n  [\"4+77=\", \"90-89=\", \"88-73=\", \"54+18=\", \"23-19=\"],\n  [\"91-48=\", \"49-24=\", \"74+12=\", \"29+15=\", \"33+37=\"],\n  [\"42-3=\", \"17+56=\", \"64+27=\", \"29+17=\", \"59-49=\"],\n  [\"17+50=\", \"57-48=\", \"82-56=\", \"86-50=\", \"63+21=\"],\n];\n\n// 1. Update the title paragraph (the date) - it is the first paragraph\n//    in the body, before the table.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\nconst titlePara = paragraphs.items[0];\n// Replace just the text run content, keeping paragraph/run formatting.\nconst range = titlePara.getRange();\nrange.insertText(DATE_NEW, Word.InsertLocation.replace);\n\n// 2. Update every cell of the (single) table with the new values.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\ntable.values = TABLE_NEW;\n\nawait context.sync();\n", "ps1": "# Update the date line and every math-problem cell in the table to the\n# new values from the commit. Every \"old\" value in the document is\n# unique, so a plain Find/Replace-All per pair is unambiguous and the\n# replacement order does not matter (no \"new\" value collides with any\n# \"old\" value).\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @('2025-10-04 Saturday', '2025-10-05 Sunday'),\n    @('21-2=', '34+23='),\n    @('27+28=', '55-35='),\n    @('84-44=', '16+31='),\n    @('14+74=', '52-4='),\n    @('26-8=', '99-34='),\n    @('86-51=', '33-1='),\n    @('16+30=', '7+21='),\n    @('43-2=', '79-35='),\n    @('42+51=', '53+30='),\n    @('62+10=', '31-19='),\n    @('30-7=', '1+35='),\n    @('68+30=', '73-17='),\n    @('42+26=', '97-72='),\n    @('89-74=', '46-22='),\n    @('89-59=', '23-6='),\n    @('50-30=', '2+35='),\n    @('61+34=', '45+40='),\n    @('88-45=', '25-13='),\n    @('36-30=', '67-25='),\n    @('48+29=', '43+33='),\n    @('39-1=', '2+40='),\n    @('41+56=', '67+9='),\n    @('81-59=', '68-34='),\n    @('10+73=', '46-16='),\n    @('45-5=', '55-6='),\n    @('75-58=', '43+25='),\n    @('58+40=', '56-9='),\n    @('49+3=', '69+27='),\n    @('52+10=', '99-7='),\n    @('6+35=', '51+23='),\n    @('51+4=', '45-4='),\n    @('62+5=', '1+46='),\n    @('85-64=', '30+33='),\n    @('44+40=', '99-73='),\n    @('94-28=', '3+60='),\n    @('37-11=', '89-19='),\n    @('65+29=', '5+43='),\n    @('46+47=', '97-47='),\n    @('55-14=', '26-23='),\n    @('54-33=', '86-22='),\n    @('52-0=', '4+54='),\n    @('45-44=', '51-5='),\n    @('45+10=', '81-81='),\n    @('28+15=', '87-21='),\n    @('40-12=', '92-70='),\n    @('43+7=', '88-27='),\n    @('79-6=', '12+38='),\n    @('96-52=', '93-55='),\n    @('11+73=', '86-20='),\n    @('83-51=', '92-15='),\n    @('94-86=', '41+35='),\n    @('30-9=', '97-42='),\n    @('80-17=', '53-16='),\n    @('90-9=', '31-8='),\n    @('89+7=', '10+82='),\n    @('59+33=', '47-23='),\n    @('58+27=', '19-15='),\n    @('97-43=', '84-56='),\n    @('52+45=', '53-9='),\n    @('70-43=', '15+42='),\n    @('33-12=', '72-13='),\n    @('78-67=', '98-1='),\n    @('87+0=', '78+17='),\n    @('58-7=', '92-79='),\n    @('91-65=', '18+70='),\n    @('84-75=', '81-25='),\n    @('93-36=', '73-20='),\n    @('89-73=', '21+77='),\n    @('0+19=', '11+36='),\n    @('19+15=', '81-47='),\n    @('48+42=', '78+18='),\n    @('43-40=', '46-9='),\n    @('47+5=', '30+42='),\n    @('9+13=', '3+72='),\n    @('79-28=', '43+46='),\n    @('49-2=', '7+56='),\n    @('46+7=', '35-7='),\n    @('42+54=', '24+67='),\n    @('8+31=', '92-89='),\n    @('53-17=', '40+50='),\n    @('93-52=', '4+77='),\n    @('81-54=', '90-89='),\n    @('25+27=', '88-73='),\n    @('78-0=', '54+18='),\n    @('90-32=', '23-19='),\n    @('56-53=', '91-48='),\n    @('11+56=', '49-24='),\n    @('84-27=', '74+12='),\n    @('53-13=', '29+15='),\n    @('96-65=', '33+37='),\n    @('16+61=', '42-3='),\n    @('97-10=', '17+56='),\n    @('49+7=', '64+27='),\n    @('3+81=', '29+17='),\n    @('54-3=', '59-49='),\n    @('95-85=', '17+50='),\n    @('45+42=', '57-48='),\n    @('38-19=', '82-56='),\n    @('32+19=', '86-50='),\n    @('20+1=', '63+21=')\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    # wdFindContinue = 1, wdReplaceAll = 2\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null\n}\n"}
